$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "mappa aggiornata con nuovi concetti" - update the "Positive" row (row 3)
# values for the three concept columns (Gastro-intestinal, Cardio, Thyroid)
$ws.Range("B3").Value = 80
$ws.Range("C3").Value = 75
$ws.Range("D3").Value = 90

# Reflect the new selection spanning the updated cells (B3:D3)
$ws.Range("B3:D3").Select()
